$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Inheritance tree diagram" screenshot requirement to be specific
# to the Creature inheritance diagram, and record its filename.
$ws.Range("B2").Value = "Creature Inheritance tree diagram"
$ws.Range("C2").Value = "Creature Inheritance Drawing.png"

# Highlight the newly-added row's status cell with a green fill (Theme color
# "Green, Accent 6").
$ws.Range("A2").Interior.ThemeColor = 10

# Move the active selection, matching the saved state of the workbook.
$ws.Range("H12").Select()

$wb.Save()
